$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label "NO.CLIENTE" -> "NO.CLIENTE:"
$ws.Range("A6").Value = "NO.CLIENTE:"

# Clear out the sample data row (values only, keep formatting)
$ws.Range("A10:E10").ClearContents()

# Column A grew a bit wider to fit the new label
$ws.Columns("A").ColumnWidth = 15.29

# Labels above (AGENTE/CLIENTE) switch from centered to left aligned
$ws.Range("B3:B4").HorizontalAlignment = -4131

# Cursor ends up on D6 when the file was saved
$ws.Range("D6").Select() | Out-Null
